$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted at row 90 of this "Vega Modelo de
# Temuco - Camote" sub-table; every pre-existing row from 90 down to 112
# shifts down by one (to 91..113). Insert a row at 90 so Excel shifts the
# rest of the table down and extends the sheet dimension to R113.
$ws.Rows.Item(90).Insert()

# Populate the newly inserted row 90 with the new record's data.
$ws.Range("A90").Value = 10
$ws.Range("B90").Value = "Vega Modelo de Temuco"
$ws.Range("C90").Value = "La Araucanía"
$ws.Range("D90").Value = "2022-08-12"
$ws.Range("E90").Value = 9
$ws.Range("F90").Value = 100114002
$ws.Range("G90").Value = "Camote"
$ws.Range("H90").Value = "Sin especificar"
$ws.Range("I90").Value = "Primera"
$ws.Range("J90").Value = 30
$ws.Range("K90").Value = 20000
$ws.Range("L90").Value = 20000
$ws.Range("M90").Value = 20000
$ws.Range("N90").Value = "$/malla 20 kilos"
$ws.Range("O90").Value = "Perú"
$ws.Range("P90").Value = 1000
$ws.Range("Q90").Value = 20
$ws.Range("R90").Value = "Hortaliza"
